$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("41_deg_from_july_13")
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Shift the existing data (old rows 18-69) down by 8 rows to rows 26-77,
# making room for the new "30 deg" block (8 new rows).
$src = $ws.Range("A18:B69")
$dest = $ws.Range("A26:B77")
$src.Copy($dest)

# New column-A values/style for the inserted rows (30.25567, vertically
# centered Arial style matching the new cellXf entry).
$colA = $ws.Range("A18:A25")
$ws2.Range("A74").Copy()
$colA.PasteSpecial(-4122)
$colA.VerticalAlignment = -4108

# New column-B values use the plain/default style (no explicit style id).
$colB = $ws.Range("B18:B25")
$ws1.Range("A2").Copy()
$colB.PasteSpecial(-4122)

# These new rows use the sheet's default row height/format (no thick
# bottom border, no custom height) unlike the thick-bordered rows below.
$newRows = $ws.Range("A18:B25").EntireRow
$newRows.AutoFit()

$ws.Range("A18").Value = 30.25567
$ws.Range("A19").Value = 30.25567
$ws.Range("A20").Value = 30.25567
$ws.Range("A21").Value = 30.25567
$ws.Range("A22").Value = 30.25567
$ws.Range("A23").Value = 30.25567
$ws.Range("A24").Value = 30.25567
$ws.Range("A25").Value = 30.25567

$ws.Range("B18").Value = 6.720393
$ws.Range("B19").Value = 9.103709
$ws.Range("B20").Value = 9.229623
$ws.Range("B21").Value = 9.092536
$ws.Range("B22").Value = 9.117204
$ws.Range("B23").Value = 8.920191
$ws.Range("B24").Value = 9.085137
$ws.Range("B25").Value = 6.76458

# Update the sheet view to match where the user left the selection/scroll.
$ws.Range("E19").Select()
$excel.ActiveWindow.ScrollRow = 8

Write-Output "done"
